# Applies the periodic cryptos-list price/volume refresh described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '30.817.20'
Set-TextValue 'E2' '  +2.61%  '
Set-TextValue 'D3' '1.899.88'
Set-TextValue 'E3' '  +0.89%  '
Set-TextValue 'D4' '1.002'
Set-TextValue 'E4' '  +0.15%  '
Set-TextValue 'D5' '246.08'
Set-TextValue 'E5' '  +1.45%  '
Set-TextValue 'E6' '  +0.14%  '
Set-TextValue 'D7' '0.4973'
Set-TextValue 'E7' '  +0.31%  '
Set-TextValue 'D8' '0.2975'
Set-TextValue 'E8' '  +1.68%  '
Set-TextValue 'E9' '  +3.49%  '
Set-TextValue 'D10' '1.906.59'
Set-TextValue 'E10' '  +1.39%  '
Set-TextValue 'D11' '17.23'
Set-TextValue 'E11' '  +2.89%  '
Set-TextValue 'D12' '0.07329'
Set-TextValue 'E12' '  +2.21%  '
Set-TextValue 'D13' '91.59'
Set-TextValue 'E13' '  +6.87%  '
Set-TextValue 'D14' '5.095'
Set-TextValue 'E14' '  +5.34%  '
Set-TextValue 'D15' '0.6792'
Set-TextValue 'E15' '  +2.56%  '
Set-TextValue 'D16' '30.821.29'
Set-TextValue 'E16' '  +2.73%  '
Set-TextValue 'D17' '0.000008052'
Set-TextValue 'D18' '13.42'
Set-TextValue 'D19' '0.9992'
Set-TextValue 'E19' '  -0.10%  '
Set-TextValue 'D20' '2.154.55'
Set-TextValue 'E20' '  +1.52%  '
Set-TextValue 'D21' '1.003'
Set-TextValue 'E21' '  +0.13%  '
Set-TextValue 'D22' '4.879'
Set-TextValue 'E22' '  +2.74%  '
Set-TextValue 'D23' '180.51'
Set-TextValue 'E23' '  +34.02%  '
Set-TextValue 'D24' '6.104'
Set-TextValue 'E24' '  +9.23%  '
Set-TextValue 'D25' '9.355'
Set-TextValue 'E25' '  +3.03%  '
Set-TextValue 'D26' '155.32'
Set-TextValue 'E26' '  +3.36%  '
Set-TextValue 'E27' '  +11.99%  '
Set-TextValue 'D28' '1.946'
Set-TextValue 'E28' '  +1.85%  '
Set-TextValue 'D29' '1.392'
Set-TextValue 'E29' '  +1.55%  '
Set-TextValue 'D30' '4.367'
Set-TextValue 'E30' '  +5.09%  '
Set-TextValue 'D31' '0.08962'
Set-TextValue 'E31' '  +3.35%  '
Set-TextValue 'D32' '4.054'
Set-TextValue 'E32' '  +3.14%  '
Set-TextValue 'D33' '0.05296'
Set-TextValue 'E33' '  +6.22%  '
Set-TextValue 'D34' '0.7533'
Set-TextValue 'E34' '  +6.71%  '
Set-TextValue 'D35' '1.145'
Set-TextValue 'E35' '  +4.01%  '
Set-TextValue 'D36' '2.698'
Set-TextValue 'E36' '  +1.63%  '
Set-TextValue 'D37' '0.01908'
Set-TextValue 'E37' '  +12.35%  '
Set-TextValue 'D38' '2.710'
Set-TextValue 'E38' '  +0.26%  '
Set-TextValue 'E39' '  -0.03%  '
Set-TextValue 'D40' '0.9374'
Set-TextValue 'E40' '  +0.94%  '
Set-TextValue 'E41' '  +4.93%  '
Set-TextValue 'D42' '105.98'
Set-TextValue 'E42' '  +3.82%  '
Set-TextValue 'D43' '5.864'
Set-TextValue 'E43' '  -1.70%  '
Set-TextValue 'D44' '1.002'
Set-TextValue 'E44' '  +0.12%  '
Set-TextValue 'D45' '7.742'
Set-TextValue 'E45' '  +4.10%  '
Set-TextValue 'D46' '0.1377'
Set-TextValue 'E46' '  +9.61%  '
Set-TextValue 'D47' '0.05845'
Set-TextValue 'E47' '  +2.98%  '
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue 'D48' '33.60'
Set-TextValue 'E48' '  +3.42%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D49' '0.3919'
Set-TextValue 'E49' '  +5.80%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D50' '8.583'
Set-TextValue 'E50' '  +4.49%  '
Set-TextValue 'D51' '1.393'
Set-TextValue 'E51' '  +4.13%  '
